$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("coisas a fazer")

# Row 11: add index number and new task text
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Fazer a divisão do banco de dados de testes e de produção"

# Row 12: add index number and new task text
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Deixar o projeto online com base de dados ativa"

# Copy style from row 10 (A10/B10) to the newly filled cells in rows 11 and 12
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11:A12").PasteSpecial(-4122) | Out-Null

$ws.Range("B10").Copy() | Out-Null
$ws.Range("B11:B12").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Update the selection to B13, matching the diff
$ws.Activate()
$ws.Range("B13").Select() | Out-Null
